# Updated cryptos list on Thu Oct 12 09:36:45 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row from the
# latest coinranking.com snapshot, and swaps rank #50 from Algorand to
# BabyDogeCoin (name, link, price, volume) to mirror the live top-50 list.
#
# Price values that are plain decimals (e.g. "21.46") would otherwise be
# auto-converted to numbers by the COM Value setter (losing the fixed-format
# text such as trailing zeros, e.g. "3.70" -> 3.7). Briefly marking the cell
# as Text (NumberFormat "@") before the write keeps the literal string, and
# resetting the style back to Normal afterwards leaves formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.793.03"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "1.554.86"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "1.558.54"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "26.801.97"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "1.367.57"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.922"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.804"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "1.692.79"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0513"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.99%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0979"
$ws.Range("E51").Value = "  -0.17%  "
